# Zero out all the numeric values in the "Valores" sheet data range (A1:L16).
# This mirrors the commit's behaviour of reporting all zeros when the
# competencia/ano could not be matched (Excel "apresentacao" bug workaround).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valores")

$rows = 16
$cols = 12

$values = New-Object 'object[,]' $rows, $cols
for ($r = 0; $r -lt $rows; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $values[$r, $c] = 0
    }
}

$range = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($rows, $cols))
$range.Value = $values
